$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "https://department-of-veterans-affairs.github.io/mhv-fhir-phr-mapping/ValueSet/ConditionVerificationVS"
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-08-22T16:36:15-05:00"
$ws.Range("B9").Value = "VA Digital Services"
